# feat: add 2022-Q4 data
#
# 1) Insert a new worksheet "2022-Q4" right after "总计" (the summary sheet),
#    populated with the quarter's fund-holding rows.
# 2) Prepend a corresponding row to the "总计" summary sheet, shifting the
#    existing quarters down by one row.

$wb = $excel.ActiveWorkbook
$summary = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. New "2022-Q4" sheet, inserted right after the summary sheet.
# ---------------------------------------------------------------------------
$q4 = $wb.Worksheets.Add($null, $summary)
$q4.Name = "2022-Q4"

# Header row.
$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

# Match the header style used throughout the workbook (bold, bordered,
# centered) by copying it from an existing quarter sheet's header row.
$styleSource = $wb.Worksheets.Item("2022-Q2")
$styleSource.Range("B1:H1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Columns B:G hold fund codes / size / position numbers that must stay text
# (leading zeros in fund codes, fixed decimal display) - force text format
# before writing so Excel doesn't coerce them into numbers.
$q4.Range("B2:G5").NumberFormat = "@"

$q4Rows = @(
    @("013051", "汇泉臻心致远混合A", "2.02", "76.67", "3.93", "0.0794", 2),
    @("013052", "汇泉臻心致远混合C", "1.37", "76.67", "3.93", "0.0538", 2),
    @("013491", "同泰金融精选股票C", "0.07", "84.75", "2.78", "0.0019", 10),
    @("013490", "同泰金融精选股票A", "0.03", "84.75", "2.78", "0.0008", 10)
)

for ($i = 0; $i -lt $q4Rows.Length; $i++) {
    $r = $i + 2
    $row = $q4Rows[$i]
    $q4.Range("A$r").Value = $i
    $q4.Range("B$r").Value = $row[0]
    $q4.Range("C$r").Value = $row[1]
    $q4.Range("D$r").Value = $row[2]
    $q4.Range("E$r").Value = $row[3]
    $q4.Range("F$r").Value = $row[4]
    $q4.Range("G$r").Value = $row[5]
    $q4.Range("H$r").Value = $row[6]
}

# Column A's index cells carry the bordered/centered style too - copy it
# from the header (same visual style used across all quarter sheets).
$q4.Range("B1").Copy()
$q4.Range("A2:A5").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2. Prepend the "2022-Q4" row to the "总计" summary sheet.
# ---------------------------------------------------------------------------
$summaryRows = @(
    @("2022-Q4", 4, 0.14),
    @("2022-Q2", 5, 0.97),
    @("2022-Q1", 1, 0.08),
    @("2021-Q4", 4, 1.62),
    @("2021-Q3", 16, 1.59),
    @("2021-Q2", 19, 0.85),
    @("2021-Q1", 25, 0.4),
    @("2020-Q4", 51, 1.69)
)

for ($i = 0; $i -lt $summaryRows.Length; $i++) {
    $r = $i + 2
    $row = $summaryRows[$i]
    $summary.Range("A$r").Value = $i
    $summary.Range("B$r").Value = $row[0]
    $summary.Range("C$r").Value = $row[1]
    $summary.Range("D$r").Value = $row[2]
}

# Row 9 is brand new (the table used to stop at row 8) - give its index
# cell (A9) the same style as the rest of column A.
$summary.Range("A2").Copy()
$summary.Range("A9").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
